$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6822
$ws.Range("J3").Value = 7211
$ws.Range("B4").Value = 1687
$ws.Range("H4").Value = 1711
$ws.Range("J4").Value = 1572
$ws.Range("J5").Value = 566
$ws.Range("I6").Value = 8966
$ws.Range("J6").Value = 9621
$ws.Range("B7").Value = 23320
$ws.Range("H7").Value = 26022
$ws.Range("I7").Value = 26231
$ws.Range("J7").Value = 25792

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J7").Value = 746
$ws.Range("J8").Value = 1621
$ws.Range("J10").Value = 189
$ws.Range("J11").Value = 444
$ws.Range("J14").Value = 137
$ws.Range("J15").Value = 306
$ws.Range("J19").Value = 753
$ws.Range("J21").Value = 71
$ws.Range("J23").Value = 238
$ws.Range("J27").Value = 152
$ws.Range("J29").Value = 1398
$ws.Range("J31").Value = 259
$ws.Range("J33").Value = 1164
$ws.Range("J34").Value = 114
$ws.Range("J36").Value = 351
$ws.Range("J37").Value = 796
$ws.Range("J40").Value = 55
$ws.Range("J42").Value = 1110
$ws.Range("J43").Value = 221
$ws.Range("J44").Value = 197
$ws.Range("J46").Value = 86
$ws.Range("J50").Value = 154
$ws.Range("J52").Value = 654
$ws.Range("J53").Value = 372
$ws.Range("J54").Value = 497
$ws.Range("J56").Value = 37
$ws.Range("J60").Value = 149
$ws.Range("B63").Value = 394
$ws.Range("H63").Value = 269
$ws.Range("I63").Value = 181
$ws.Range("J63").Value = 80
$ws.Range("J64").Value = 170
$ws.Range("J65").Value = 645
$ws.Range("J67").Value = 967
$ws.Range("J75").Value = 78
$ws.Range("J76").Value = 377
$ws.Range("J78").Value = 305
$ws.Range("J79").Value = 728
$ws.Range("J83").Value = 514
$ws.Range("J86").Value = 165
$ws.Range("J88").Value = 272
$ws.Range("J90").Value = 274
$ws.Range("J92").Value = 82
$ws.Range("J95").Value = 376
$ws.Range("J96").Value = 280
$ws.Range("J97").Value = 237
$ws.Range("J99").Value = 396
$ws.Range("B101").Value = 23320
$ws.Range("H101").Value = 26022
$ws.Range("I101").Value = 26231
$ws.Range("J101").Value = 25792

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 231
$ws.Range("J7").Value = 746

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 202
$ws.Range("J7").Value = 444

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 154
$ws.Range("J3").Value = 185
$ws.Range("J6").Value = 280
$ws.Range("J7").Value = 654

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 246
$ws.Range("J7").Value = 372

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 436
$ws.Range("J3").Value = 484
$ws.Range("J5").Value = 42
$ws.Range("J6").Value = 572
$ws.Range("J7").Value = 1621

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 153
$ws.Range("J3").Value = 192
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 514

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 388
$ws.Range("J6").Value = 405
$ws.Range("J7").Value = 1164

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 129
$ws.Range("J7").Value = 376

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 233
$ws.Range("J7").Value = 796

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 187
$ws.Range("J7").Value = 645

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 396

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 358
$ws.Range("J7").Value = 967

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 123
$ws.Range("J6").Value = 232
$ws.Range("J7").Value = 497

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 424
$ws.Range("J3").Value = 496
$ws.Range("J6").Value = 354
$ws.Range("J7").Value = 1398

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 216
$ws.Range("J6").Value = 292
$ws.Range("J7").Value = 753

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 64
$ws.Range("J7").Value = 377

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 237
$ws.Range("J4").Value = 46
$ws.Range("J6").Value = 586
$ws.Range("J7").Value = 1110

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 33
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 81
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 203
$ws.Range("J3").Value = 245
$ws.Range("J6").Value = 216
$ws.Range("J7").Value = 728

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 47
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 113
$ws.Range("J7").Value = 351

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 65
$ws.Range("J6").Value = 138
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 87
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 274

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 132
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 55
